$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 828, shifting existing rows 828:869 down to 829:870
$ws.Rows("828:828").Insert()

# The date column (A) must stay a plain text value like "2026/02/18" rather than
# being auto-converted to a date serial number, so force a text number format
# before assigning, then restore the default "Normal" style/format afterwards.
$ws.Range("A828").NumberFormat = "@"
$ws.Range("A828").Value = "2026/02/18"
$ws.Range("A828").Style = "Normal"

$ws.Range("B828").Value = "水"
$ws.Range("C828").Value = 18
$ws.Range("D828").Value = 201
